# "Generate Report for Handoff"
#
# A fresh handoff run completed for the four files that were previously
# sitting at "low" priority / queued for handoff
# (0b617f20-*, 97ff1396-*, a01547b3-*, faca6966-*). For each locale sheet
# ("zh-cn" and "de-de") the report generator:
#   - raises the Priority (column E) for those rows from "low" to "ht"
#   - stamps the Latest Handoff Datetime (column H) for those rows with the
#     timestamp of the new handoff xliff generation
#
# The "Overview" sheet's "Latest HO Xliff Generate Date" column mirrors the
# most recent per-locale handoff timestamp for each file, so it is refreshed
# with the same new date (de-de ran last, so that's the value that shows
# through there).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhHandoffTime = "2016-09-05 00:36:26"
$deHandoffTime = "2016-09-05 00:36:31"

foreach ($row in 4..7) {
    $zhcn.Range("E$row").Value = "ht"
    $zhcn.Range("H$row").Value = $zhHandoffTime

    $dede.Range("E$row").Value = "ht"
    $dede.Range("H$row").Value = $deHandoffTime

    $overview.Range("G$row").Value = $deHandoffTime
}
